$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sprint 3": a new "S" day column was inserted before the old column
# E ("M"), shifting the D:Q day columns one to the right (the trailing
# Saturday that fell off the end of the Q column is simply dropped), and a
# new backlog row ("Desgin Logo") of work was added as row 5.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sprint 3")

# Header row (row 1): D1:Q1 day-of-week labels after the shift.
$headerVals = @("S","S","M","T","W","T","F","S","S","M","T","W","T","F")
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $col = 4 + $i   # D = column 4
    $ws.Cells.Item(1, $col).Value2 = $headerVals[$i]
}

# Row 2 ("Display search results"): values shift right from G2 into H2:J2.
$ws.Cells.Item(2, 8).Value2 = 8   # H2
$ws.Cells.Item(2, 9).Value2 = 8   # I2
$ws.Cells.Item(2, 10).Value2 = 8  # J2

# Row 3 ("Search Activity"): G3 gains the old F3 value (5) and the old G3
# value (2) shifts right into H3, with two more 2's following into I3/J3.
$ws.Cells.Item(3, 7).Value2 = 5   # G3
$ws.Cells.Item(3, 8).Value2 = 2   # H3
$ws.Cells.Item(3, 9).Value2 = 2   # I3
$ws.Cells.Item(3, 10).Value2 = 2  # J3

# Row 4 ("App background service"): values shift right from G4 into H4:J4.
$ws.Cells.Item(4, 8).Value2 = 5   # H4
$ws.Cells.Item(4, 9).Value2 = 5   # I4
$ws.Cells.Item(4, 10).Value2 = 5  # J4

# Row 5: brand-new backlog item "Desgin Logo".
$ws.Cells.Item(5, 1).Value2 = 10              # A5 - Product ID
$ws.Cells.Item(5, 2).Value2 = "Desgin Logo"   # B5 - Task
$ws.Cells.Item(5, 3).Value2 = 3               # C5 - Start hours
$ws.Cells.Item(5, 4).Value2 = 3               # D5
$ws.Cells.Item(5, 5).Value2 = 3               # E5
$ws.Cells.Item(5, 6).Value2 = 3               # F5
$ws.Cells.Item(5, 7).Value2 = 3               # G5
$ws.Cells.Item(5, 8).Value2 = 3               # H5
$ws.Cells.Item(5, 9).Value2 = 3               # I5
$ws.Cells.Item(5, 10).Value2 = 3              # J5

# ---------------------------------------------------------------------------
# Sheet "Product": only the active selection changed (D11 -> J18).
# ---------------------------------------------------------------------------
$wsProduct = $wb.Worksheets.Item("Product")
$wsProduct.Range("J18").Select()

# ---------------------------------------------------------------------------
# "Sprint 3" stays the active/selected tab, now with R8 selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("R8").Select()
